# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback datetime
# timestamps to reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: G2 - "Latest HO Xliff Generate Date" for 130d674f-...md
$wsOverview.Range("G2").Value = "2016-10-19 11:33:43"

# zh-cn: H2 - Correspond Handoff Datetime, K2 - Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-10-19 11:33:32"
$wsZhCn.Range("K2").Value = "2016-10-19 11:34:12"

# de-de: H2 - Correspond Handoff Datetime (mirrors Overview G2's value)
$wsDeDe.Range("H2").Value = "2016-10-19 11:33:43"
# de-de: K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-10-19 11:34:29"
